# Rapise 6.6 note update:
#  - Add a second column (B) next to the existing dropdown values in column A.
#  - Column B duplicates column A's list, except the first (header) row which
#    uses a slightly different key string ("CrmChangeArea.name" instead of
#    "Crm.ChangeArea.name").
#  - Give column B a custom width and move the active selection to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B values (mirrors column A, row for row) -----------------------
$ws.Range("B1").Value = "CrmChangeArea.name"
$ws.Range("B2").Value = "Sales"
$ws.Range("B3").Value = "App Settings"
$ws.Range("B4").Value = "Sales Insights settings"
$ws.Range("B5").Value = "Help and Support"

# Header cell B1 gets the same bold style as A1.
$ws.Range("B1").Font.Bold = $true

# --- Column B width ----------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 21.25

# --- Move / update the active selection to B8 -------------------------------
$ws.Range("B8").Select() | Out-Null
